# Scheduled market-price refresh: update the price/profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) for the rows whose quoted prices changed, across all
# eight item-category sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2607
$ws.Range("I15").Value = 2607
$ws.Range("K15").Value = 7821
$ws.Range("M15").Value = -7652
$ws.Range("H33").Value = 683.1739
$ws.Range("I33").Value = 504.2143
$ws.Range("J33").Value = 961.55554
$ws.Range("K33").Value = 504.2143
$ws.Range("L33").Value = 961.55554
$ws.Range("M33").Value = -275.2143
$ws.Range("N33").Value = -1419.55554
$ws.Range("H86").Value = 21480.6
$ws.Range("I86").Value = 21480.6
$ws.Range("K86").Value = 21480.6
$ws.Range("M86").Value = -20357.6
$ws.Range("H89").Value = 21480.6
$ws.Range("I89").Value = 21480.6
$ws.Range("K89").Value = 107403
$ws.Range("M89").Value = -101787
$ws.Range("H135").Value = 1418.4166
$ws.Range("I135").Value = 1025.3
$ws.Range("K135").Value = 9227.699999999999
$ws.Range("M135").Value = -6692.699999999999
$ws.Range("H137").Value = 3919.5
$ws.Range("I137").Value = 3821.7144
$ws.Range("J137").Value = 4212.857
$ws.Range("K137").Value = 11465.1432
$ws.Range("L137").Value = 12638.571
$ws.Range("M137").Value = -8915.143199999999
$ws.Range("N137").Value = -17738.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 597.0526
$ws.Range("I74").Value = 569.94446
$ws.Range("J74").Value = 1085
$ws.Range("K74").Value = 569.94446
$ws.Range("L74").Value = 1085
$ws.Range("M74").Value = 304.05554
$ws.Range("N74").Value = -2833
$ws.Range("H77").Value = 597.0526
$ws.Range("I77").Value = 569.94446
$ws.Range("J77").Value = 1085
$ws.Range("K77").Value = 2849.7223
$ws.Range("L77").Value = 5425
$ws.Range("M77").Value = 1518.2777
$ws.Range("N77").Value = -14161
$ws.Range("H88").Value = 2163.25
$ws.Range("I88").Value = 1901.5
$ws.Range("J88").Value = 2425
$ws.Range("K88").Value = 1901.5
$ws.Range("L88").Value = 2425
$ws.Range("M88").Value = -1495.5
$ws.Range("N88").Value = -3237
$ws.Range("H91").Value = 2163.25
$ws.Range("I91").Value = 1901.5
$ws.Range("J91").Value = 2425
$ws.Range("K91").Value = 1901.5
$ws.Range("L91").Value = 2425
$ws.Range("M91").Value = -497.5
$ws.Range("N91").Value = -5233

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7016.077
$ws.Range("I20").Value = 7130.9
$ws.Range("J20").Value = 6633.3335
$ws.Range("K20").Value = 7130.9
$ws.Range("L20").Value = 6633.3335
$ws.Range("M20").Value = -6883.9
$ws.Range("N20").Value = -7127.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3131.68
$ws.Range("I31").Value = 2458.4062
$ws.Range("J31").Value = 4328.6113
$ws.Range("K31").Value = 2458.4062
$ws.Range("L31").Value = 4328.6113
$ws.Range("M31").Value = -2163.4062
$ws.Range("N31").Value = -4918.6113
$ws.Range("H34").Value = 3131.68
$ws.Range("I34").Value = 2458.4062
$ws.Range("J34").Value = 4328.6113
$ws.Range("K34").Value = 2458.4062
$ws.Range("L34").Value = 4328.6113
$ws.Range("M34").Value = -2256.4062
$ws.Range("N34").Value = -4732.6113
$ws.Range("H51").Value = 11924.75
$ws.Range("J51").Value = 11924.75
$ws.Range("L51").Value = 11924.75
$ws.Range("N51").Value = -13396.75
$ws.Range("H61").Value = 11924.75
$ws.Range("J61").Value = 11924.75
$ws.Range("L61").Value = 11924.75
$ws.Range("N61").Value = -12620.75
$ws.Range("H111").Value = 39800
$ws.Range("J111").Value = 39800
$ws.Range("L111").Value = 39800
$ws.Range("N111").Value = -47980
$ws.Range("H134").Value = 2041
$ws.Range("I134").Value = 968.05554
$ws.Range("J134").Value = 4800
$ws.Range("K134").Value = 2904.16662
$ws.Range("L134").Value = 14400
$ws.Range("M134").Value = -369.16662
$ws.Range("N134").Value = -19470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 41.53846
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 52.2
$ws.Range("K2").Value = 36
$ws.Range("L2").Value = 313.2
$ws.Range("M2").Value = 77
$ws.Range("N2").Value = -539.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4934.5454
$ws.Range("I70").Value = 4912.5
$ws.Range("J70").Value = 4993.3335
$ws.Range("K70").Value = 4912.5
$ws.Range("L70").Value = 4993.3335
$ws.Range("M70").Value = -4642.5
$ws.Range("N70").Value = -5533.3335
$ws.Range("H73").Value = 4934.5454
$ws.Range("I73").Value = 4912.5
$ws.Range("J73").Value = 4993.3335
$ws.Range("K73").Value = 4912.5
$ws.Range("L73").Value = 4993.3335
$ws.Range("M73").Value = -3976.5
$ws.Range("N73").Value = -6865.3335
$ws.Range("H80").Value = 3839.4211
$ws.Range("I80").Value = 3912.5386
$ws.Range("J80").Value = 3681
$ws.Range("K80").Value = 3912.5386
$ws.Range("L80").Value = 3681
$ws.Range("M80").Value = -2914.5386
$ws.Range("N80").Value = -5677
$ws.Range("H83").Value = 3839.4211
$ws.Range("I83").Value = 3912.5386
$ws.Range("J83").Value = 3681
$ws.Range("K83").Value = 19562.693
$ws.Range("L83").Value = 18405
$ws.Range("M83").Value = -14570.693
$ws.Range("N83").Value = -28389
$ws.Range("H97").Value = 1955.4546
$ws.Range("I97").Value = 1001.4286
$ws.Range("J97").Value = 3625
$ws.Range("K97").Value = 1001.4286
$ws.Range("L97").Value = 3625
$ws.Range("M97").Value = -505.4286
$ws.Range("N97").Value = -4617

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 57504.5
$ws.Range("J18").Value = 57504.5
$ws.Range("L18").Value = 57504.5
$ws.Range("N18").Value = -57848.5
$ws.Range("H93").Value = 6450
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 6450
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 6450
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -8946
$ws.Range("H132").Value = 3166.8276
$ws.Range("I132").Value = 1937.091
$ws.Range("J132").Value = 3918.3333
$ws.Range("K132").Value = 5811.272999999999
$ws.Range("L132").Value = 11754.9999
$ws.Range("M132").Value = -3281.272999999999
$ws.Range("N132").Value = -16814.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3848810.5
$ws.Range("I126").Value = 1953.2
$ws.Range("J126").Value = 16671668
$ws.Range("K126").Value = 5859.6
$ws.Range("L126").Value = 50015004
$ws.Range("M126").Value = -3389.6
$ws.Range("N126").Value = -50019944
